$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set every Fitness value (column C, rows 2-252) to 7639
$ws.Range("C2:C252").Value = 7639
